$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Rufino Imóveis"
$ws.Range("B5").Value = " R. Treze de Maio, 33 - Centro"
$ws.Range("C5").Value = -21.38934490323197
$ws.Range("D5").Value = -42.69477177427834
$ws.Range("G5").Value = "https://i.imgur.com/HU7Ezto.png"
